# Updates betting odds values on Sheet1 as per the source diff.
# Cell values are updated directly via the Excel COM object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.5
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 13
$ws.Range("AE2").Value = 23
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 19
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 10
$ws.Range("AQ2").Value = 34
$ws.Range("AU2").Value = 10
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 34
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 201
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 2.5
$ws.Range("L3").Value = 5.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 7.5
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 17
$ws.Range("AB3").Value = 41
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 51
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 10
$ws.Range("AQ3").Value = 34
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 29
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 151
$ws.Range("L6").Value = 3.6
$ws.Range("W6").Value = 7.6
$ws.Range("X6").Value = 10.75
$ws.Range("AA6").Value = 18
$ws.Range("AB6").Value = 28
$ws.Range("AD6").Value = 6.2
$ws.Range("AG6").Value = 9.25
$ws.Range("AH6").Value = 16.5
$ws.Range("AK6").Value = 28
$ws.Range("AL6").Value = 35
$ws.Range("AN6").Value = 4.1
$ws.Range("AP6").Value = 18
$ws.Range("AQ6").Value = 45
$ws.Range("AR6").Value = 70
$ws.Range("AT6").Value = 2.57
$ws.Range("AU6").Value = 6.7
$ws.Range("G8").Value = 3.85
$ws.Range("H8").Value = 3.4
$ws.Range("J8").Value = 4.3
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 2.42
$ws.Range("P8").Value = 3.15
$ws.Range("R8").Value = 1.88
$ws.Range("W8").Value = 11.25
$ws.Range("Y8").Value = 13
$ws.Range("AB8").Value = 40
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 6.6
$ws.Range("AG8").Value = 7.7
$ws.Range("AH8").Value = 9.5
$ws.Range("AJ8").Value = 16.5
$ws.Range("AK8").Value = 14
$ws.Range("AN8").Value = 5.7
$ws.Range("AO8").Value = 22
$ws.Range("AP8").Value = 27
$ws.Range("AQ8").Value = 120
$ws.Range("AY8").Value = 17
$ws.Range("AZ8").Value = 32
$ws.Range("BA8").Value = 60
$ws.Range("G17").Value = 2.1
$ws.Range("I17").Value = 3.25
$ws.Range("L17").Value = 3.6
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("AC17").Value = 13
$ws.Range("AH17").Value = 19
$ws.Range("AL17").Value = 29
$ws.Range("AN17").Value = 4.33
$ws.Range("AQ17").Value = 34
$ws.Range("AU17").Value = 7.5
$ws.Range("H25").Value = 3.45
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 2.67
$ws.Range("P25").Value = 3.65
$ws.Range("R25").Value = 2.02
$ws.Range("S25").Value = 1.35
$ws.Range("T25").Value = 2.95
$ws.Range("X25").Value = 11.25
$ws.Range("AB25").Value = 23
$ws.Range("AD25").Value = 6.9
$ws.Range("AE25").Value = 13
$ws.Range("AG25").Value = 11
$ws.Range("AH25").Value = 17
$ws.Range("AP25").Value = 17.5
$ws.Range("AT25").Value = 2.95
$ws.Range("AW25").Value = 5.1
$ws.Range("N26").Value = 5.7
